$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values per row (destination row -> D, J, K, L, M, P), derived from
# shuffling the existing (Fecha, Volumen, Precio mínimo, Precio máximo,
# Precio promedio ponderado, Precio $/Kg) tuples among rows 2-16.
$data = @{
  2  = @(44365, 55, 5000, 5000, 5000, 5000)
  3  = @(44313, 20, 4000, 4000, 4000, 4000)
  4  = @(44316, 20, 4000, 4000, 4000, 4000)
  5  = @(44280, 55, 4000, 4000, 4000, 4000)
  6  = @(44497, 20, 4000, 4000, 4000, 4000)
  7  = @(44291, 35, 4000, 4000, 4000, 4000)
  8  = @(44498, 40, 4000, 4000, 4000, 4000)
  9  = @(44509, 20, 4000, 4000, 4000, 4000)
  10 = @(44259, 30, 4000, 4000, 4000, 4000)
  11 = @(44176, 10, 4000, 4000, 4000, 4000)
  12 = @(44301, 40, 3000, 3000, 3000, 3000)
  13 = @(44504, 55, 4000, 4000, 4000, 4000)
  14 = @(44312, 50, 4000, 4000, 4000, 4000)
  15 = @(44315, 40, 4000, 4000, 4000, 4000)
  16 = @(44508, 30, 4000, 4000, 4000, 4000)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # D: Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]  # J: Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]  # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[3]  # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[4]  # M: Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals[5]  # P: Precio $/Kg
}
